$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A34").Value = "GRT-USD"
